# Automatic update of files.
# Column C ("Förändrad") holds the "last changed" date for every record,
# stored as an Excel date serial number. This run bumps that date by one
# day (46074 -> 46075, i.e. 2026-02-21 -> 2026-02-22) for every data row
# (rows 2 through 281).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C281").Value = 46075
